$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 125 (1-indexed), shifting existing rows 125+ down by one.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new weekly data entry.
$ws.Cells.Item(125, 1).Value = 7
$ws.Cells.Item(125, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(125, 3).Value = "Ñuble"
$ws.Cells.Item(125, 4).Value = 44452
$ws.Cells.Item(125, 5).Value = 16
$ws.Cells.Item(125, 6).Value = 100114001
$ws.Cells.Item(125, 7).Value = "Papa"
$ws.Cells.Item(125, 8).Value = "Patagonia"
$ws.Cells.Item(125, 9).Value = "1a (guarda)"
$ws.Cells.Item(125, 10).Value = 300
$ws.Cells.Item(125, 11).Value = 7000
$ws.Cells.Item(125, 12).Value = 7500
$ws.Cells.Item(125, 13).Value = 7250
$ws.Cells.Item(125, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(125, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(125, 16).Value = 290
$ws.Cells.Item(125, 17).Value = 25
$ws.Cells.Item(125, 18).Value = "Hortaliza"
